$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.255.25"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "3.320.50"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "187.45"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "585.74"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.71"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.411"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").Value = "3.907.41"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.76"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "68.509.85"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000168"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "3.326.33"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "444.77"
$ws.Range("E18").Value = "  +12.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.76"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.67"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.76"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.32"
$ws.Range("E22").Value = "  +6.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.519"
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("D25").Value = "3.482.79"
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.189"
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.28"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.994"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("E30").Value = "  +2.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.09"
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.42"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.88"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.53"
$ws.Range("E36").Value = "  +5.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.56"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.90"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.03"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.58"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.792"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.47"
$ws.Range("E42").Value = "  +2.71%  "
$ws.Range("D43").Value = "2.690.75"
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.92"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.45"
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0680"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.79"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "328.23"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0278"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.07"
$ws.Range("E50").Value = "  +5.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.993"
$ws.Range("E51").Value = "  +2.36%  "
